$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values must remain plain text even though many of them
# look like numbers (e.g. "231.67"). Excel auto-converts numeric-looking text
# to a real number on assignment, so each target cell is temporarily switched
# to the Text number format before the value is written, then its style is
# reset back to Normal so no stray formatting is left on the cell.
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "28.667.50"
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.800.38"
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "231.67"
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.5898"
$cell.Style = "Normal"
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.2768"
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06800"
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "23.27"
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07511"
$cell.Style = "Normal"
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.847.67"
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "4.762"
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.6216"
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "2.044.50"
$cell.Style = "Normal"
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.000009138"
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "75.48"
$cell.Style = "Normal"
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "28.642.84"
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "5.459"
$cell.Style = "Normal"
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "210.26"
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "11.50"
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "6.816"
$cell.Style = "Normal"
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "153.76"
$cell.Style = "Normal"
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "7.845"
$cell.Style = "Normal"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.1266"
$cell.Style = "Normal"
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "16.42"
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.417"
$cell.Style = "Normal"
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.06187"
$cell.Style = "Normal"
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.422"
$cell.Style = "Normal"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.770"
$cell.Style = "Normal"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.748"
$cell.Style = "Normal"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.729"
$cell.Style = "Normal"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.6418"
$cell.Style = "Normal"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.499"
$cell.Style = "Normal"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "2.711"
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "6.550"
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.148.03"
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.8831"
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.007"
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "100.06"
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.948.41"
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "60.28"
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.00000000111"
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.586"
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "8.356"
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.05460"
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.4481"
$cell.Style = "Normal"

# Column E ("Volume(1h)") values are percentage strings (e.g. "  -2.25%  ")
# and are never mistaken for numbers, so they can be assigned directly.
$ws.Range("E2").Value = "  -2.25%  "
$ws.Range("E3").Value = "  -1.82%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("E6").Value = "  -2.46%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  -1.33%  "
$ws.Range("E9").Value = "  -3.99%  "
$ws.Range("E10").Value = "  -1.49%  "
$ws.Range("E11").Value = "  -1.99%  "
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("E15").Value = "  -1.79%  "
$ws.Range("E16").Value = "  -8.36%  "
$ws.Range("E17").Value = "  -4.88%  "
$ws.Range("E18").Value = "  -2.25%  "
$ws.Range("E19").Value = "  -7.04%  "
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("E21").Value = "  -7.51%  "
$ws.Range("E22").Value = "  -2.06%  "
$ws.Range("E23").Value = "  -2.94%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("E26").Value = "  -2.51%  "
$ws.Range("E27").Value = "  -2.99%  "
$ws.Range("E28").Value = "  -1.28%  "
$ws.Range("E29").Value = "  -4.61%  "
$ws.Range("E30").Value = "  -1.05%  "
$ws.Range("E31").Value = "  -1.83%  "
$ws.Range("E32").Value = "  -1.87%  "
$ws.Range("E33").Value = "  -1.68%  "
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("E35").Value = "  -6.24%  "
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("E37").Value = "  -1.80%  "
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("E40").Value = "  -3.46%  "
$ws.Range("E41").Value = "  -6.24%  "
$ws.Range("E42").Value = "  -2.66%  "
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("E44").Value = "  -0.86%  "
$ws.Range("E45").Value = "  -2.06%  "
$ws.Range("E46").Value = "  -4.15%  "
$ws.Range("E47").Value = "  -3.90%  "
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("E49").Value = "  -1.91%  "
$ws.Range("E50").Value = "  -1.00%  "
$ws.Range("E51").Value = "  -1.92%  "
